$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date number format (YYYY-MM-DD HH:MM:SS, style index 2) from an existing
# date cell so new/rewritten column-D cells keep the same date styling.
$dateFmt = $ws.Range("D2").NumberFormat
$ws.Range("D260:D276").NumberFormat = $dateFmt

# Row 260
$ws.Range("A260").Value = 6
$ws.Range("B260").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C260").Value = "Metropolitana"
$ws.Range("D260").Value = 44516
$ws.Range("E260").Value = 13
$ws.Range("F260").Value = 100112052
$ws.Range("G260").Value = "Albahaca"
$ws.Range("H260").Value = "Sin especificar"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 300
$ws.Range("K260").Value = 5000
$ws.Range("L260").Value = 5500
$ws.Range("M260").Value = 5167
$ws.Range("N260").Value = "$/docena de matas"
$ws.Range("O260").Value = "Región Metropolitana"
$ws.Range("P260").Value = 861
$ws.Range("Q260").Value = 6
$ws.Range("R260").Value = "Hortaliza"

# Row 261
$ws.Range("A261").Value = 6
$ws.Range("B261").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C261").Value = "Metropolitana"
$ws.Range("D261").Value = 44516
$ws.Range("E261").Value = 13
$ws.Range("F261").Value = 100112052
$ws.Range("G261").Value = "Albahaca"
$ws.Range("H261").Value = "Sin especificar"
$ws.Range("I261").Value = "Segunda"
$ws.Range("J261").Value = 50
$ws.Range("K261").Value = 4000
$ws.Range("L261").Value = 4000
$ws.Range("M261").Value = 4000
$ws.Range("N261").Value = "$/docena de matas"
$ws.Range("O261").Value = "Región Metropolitana"
$ws.Range("P261").Value = 667
$ws.Range("Q261").Value = 6
$ws.Range("R261").Value = "Hortaliza"

# Row 262
$ws.Range("A262").Value = 6
$ws.Range("B262").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C262").Value = "Metropolitana"
$ws.Range("D262").Value = 44270
$ws.Range("E262").Value = 13
$ws.Range("F262").Value = 100112052
$ws.Range("G262").Value = "Albahaca"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 150
$ws.Range("K262").Value = 3000
$ws.Range("L262").Value = 3000
$ws.Range("M262").Value = 3000
$ws.Range("N262").Value = "$/docena de matas"
$ws.Range("O262").Value = "Región Metropolitana"
$ws.Range("P262").Value = 500
$ws.Range("Q262").Value = 6
$ws.Range("R262").Value = "Hortaliza"

# Row 263
$ws.Range("A263").Value = 6
$ws.Range("B263").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C263").Value = "Metropolitana"
$ws.Range("D263").Value = 44270
$ws.Range("E263").Value = 13
$ws.Range("F263").Value = 100112052
$ws.Range("G263").Value = "Albahaca"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Segunda"
$ws.Range("J263").Value = 90
$ws.Range("K263").Value = 2500
$ws.Range("L263").Value = 2500
$ws.Range("M263").Value = 2500
$ws.Range("N263").Value = "$/docena de matas"
$ws.Range("O263").Value = "Región Metropolitana"
$ws.Range("P263").Value = 417
$ws.Range("Q263").Value = 6
$ws.Range("R263").Value = "Hortaliza"

# Row 264
$ws.Range("A264").Value = 6
$ws.Range("B264").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C264").Value = "Metropolitana"
$ws.Range("D264").Value = 44295
$ws.Range("E264").Value = 13
$ws.Range("F264").Value = 100112052
$ws.Range("G264").Value = "Albahaca"
$ws.Range("H264").Value = "Sin especificar"
$ws.Range("I264").Value = "Primera"
$ws.Range("J264").Value = 210
$ws.Range("K264").Value = 3800
$ws.Range("L264").Value = 4000
$ws.Range("M264").Value = 3886
$ws.Range("N264").Value = "$/docena de matas"
$ws.Range("O264").Value = "Región Metropolitana"
$ws.Range("P264").Value = 648
$ws.Range("Q264").Value = 6
$ws.Range("R264").Value = "Hortaliza"

# Row 265
$ws.Range("A265").Value = 6
$ws.Range("B265").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C265").Value = "Metropolitana"
$ws.Range("D265").Value = 44217
$ws.Range("E265").Value = 13
$ws.Range("F265").Value = 100112052
$ws.Range("G265").Value = "Albahaca"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 980
$ws.Range("K265").Value = 3000
$ws.Range("L265").Value = 3500
$ws.Range("M265").Value = 3179
$ws.Range("N265").Value = "$/docena de matas"
$ws.Range("O265").Value = "Región Metropolitana"
$ws.Range("P265").Value = 530
$ws.Range("Q265").Value = 6
$ws.Range("R265").Value = "Hortaliza"

# Row 266
$ws.Range("A266").Value = 6
$ws.Range("B266").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C266").Value = "Metropolitana"
$ws.Range("D266").Value = 44217
$ws.Range("E266").Value = 13
$ws.Range("F266").Value = 100112052
$ws.Range("G266").Value = "Albahaca"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Segunda"
$ws.Range("J266").Value = 310
$ws.Range("K266").Value = 2000
$ws.Range("L266").Value = 2500
$ws.Range("M266").Value = 2258
$ws.Range("N266").Value = "$/docena de matas"
$ws.Range("O266").Value = "Región Metropolitana"
$ws.Range("P266").Value = 376
$ws.Range("Q266").Value = 6
$ws.Range("R266").Value = "Hortaliza"

# Row 267
$ws.Range("A267").Value = 6
$ws.Range("B267").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C267").Value = "Metropolitana"
$ws.Range("D267").Value = 44509
$ws.Range("E267").Value = 13
$ws.Range("F267").Value = 100112052
$ws.Range("G267").Value = "Albahaca"
$ws.Range("H267").Value = "Sin especificar"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 550
$ws.Range("K267").Value = 5000
$ws.Range("L267").Value = 6000
$ws.Range("M267").Value = 5545
$ws.Range("N267").Value = "$/docena de matas"
$ws.Range("O267").Value = "Región Metropolitana"
$ws.Range("P267").Value = 924
$ws.Range("Q267").Value = 6
$ws.Range("R267").Value = "Hortaliza"

# Row 268
$ws.Range("A268").Value = 6
$ws.Range("B268").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C268").Value = "Metropolitana"
$ws.Range("D268").Value = 44244
$ws.Range("E268").Value = 13
$ws.Range("F268").Value = 100112052
$ws.Range("G268").Value = "Albahaca"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 670
$ws.Range("K268").Value = 2500
$ws.Range("L268").Value = 3000
$ws.Range("M268").Value = 2612
$ws.Range("N268").Value = "$/docena de matas"
$ws.Range("O268").Value = "Región Metropolitana"
$ws.Range("P268").Value = 435
$ws.Range("Q268").Value = 6
$ws.Range("R268").Value = "Hortaliza"

# Row 269
$ws.Range("A269").Value = 6
$ws.Range("B269").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C269").Value = "Metropolitana"
$ws.Range("D269").Value = 44244
$ws.Range("E269").Value = 13
$ws.Range("F269").Value = 100112052
$ws.Range("G269").Value = "Albahaca"
$ws.Range("H269").Value = "Sin especificar"
$ws.Range("I269").Value = "Segunda"
$ws.Range("J269").Value = 320
$ws.Range("K269").Value = 2000
$ws.Range("L269").Value = 2000
$ws.Range("M269").Value = 2000
$ws.Range("N269").Value = "$/docena de matas"
$ws.Range("O269").Value = "Región Metropolitana"
$ws.Range("P269").Value = 333
$ws.Range("Q269").Value = 6
$ws.Range("R269").Value = "Hortaliza"

# Row 270
$ws.Range("A270").Value = 6
$ws.Range("B270").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C270").Value = "Metropolitana"
$ws.Range("D270").Value = 44307
$ws.Range("E270").Value = 13
$ws.Range("F270").Value = 100112052
$ws.Range("G270").Value = "Albahaca"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 170
$ws.Range("K270").Value = 5000
$ws.Range("L270").Value = 5000
$ws.Range("M270").Value = 5000
$ws.Range("N270").Value = "$/docena de matas"
$ws.Range("O270").Value = "Región Metropolitana"
$ws.Range("P270").Value = 833
$ws.Range("Q270").Value = 6
$ws.Range("R270").Value = "Hortaliza"

# Row 271
$ws.Range("A271").Value = 6
$ws.Range("B271").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C271").Value = "Metropolitana"
$ws.Range("D271").Value = 44273
$ws.Range("E271").Value = 13
$ws.Range("F271").Value = 100112052
$ws.Range("G271").Value = "Albahaca"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 730
$ws.Range("K271").Value = 2500
$ws.Range("L271").Value = 4000
$ws.Range("M271").Value = 3041
$ws.Range("N271").Value = "$/docena de matas"
$ws.Range("O271").Value = "Región Metropolitana"
$ws.Range("P271").Value = 507
$ws.Range("Q271").Value = 6
$ws.Range("R271").Value = "Hortaliza"

# Row 272
$ws.Range("A272").Value = 6
$ws.Range("B272").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C272").Value = "Metropolitana"
$ws.Range("D272").Value = 44273
$ws.Range("E272").Value = 13
$ws.Range("F272").Value = 100112052
$ws.Range("G272").Value = "Albahaca"
$ws.Range("H272").Value = "Sin especificar"
$ws.Range("I272").Value = "Segunda"
$ws.Range("J272").Value = 70
$ws.Range("K272").Value = 3500
$ws.Range("L272").Value = 3500
$ws.Range("M272").Value = 3500
$ws.Range("N272").Value = "$/docena de matas"
$ws.Range("O272").Value = "Región Metropolitana"
$ws.Range("P272").Value = 583
$ws.Range("Q272").Value = 6
$ws.Range("R272").Value = "Hortaliza"

# Row 273
$ws.Range("A273").Value = 6
$ws.Range("B273").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C273").Value = "Metropolitana"
$ws.Range("D273").Value = 44433
$ws.Range("E273").Value = 13
$ws.Range("F273").Value = 100112052
$ws.Range("G273").Value = "Albahaca"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 160
$ws.Range("K273").Value = 5500
$ws.Range("L273").Value = 6000
$ws.Range("M273").Value = 5781
$ws.Range("N273").Value = "$/paquete"
$ws.Range("O273").Value = "Región de Arica y Parinacota"
$ws.Range("P273").Value = 5781
$ws.Range("Q273").Value = 1
$ws.Range("R273").Value = "Hortaliza"

# Row 274
$ws.Range("A274").Value = 6
$ws.Range("B274").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C274").Value = "Metropolitana"
$ws.Range("D274").Value = 44302
$ws.Range("E274").Value = 13
$ws.Range("F274").Value = 100112052
$ws.Range("G274").Value = "Albahaca"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 190
$ws.Range("K274").Value = 6000
$ws.Range("L274").Value = 6000
$ws.Range("M274").Value = 6000
$ws.Range("N274").Value = "$/docena de matas"
$ws.Range("O274").Value = "Región Metropolitana"
$ws.Range("P274").Value = 1000
$ws.Range("Q274").Value = 6
$ws.Range("R274").Value = "Hortaliza"

# Row 275
$ws.Range("A275").Value = 6
$ws.Range("B275").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C275").Value = "Metropolitana"
$ws.Range("D275").Value = 44179
$ws.Range("E275").Value = 13
$ws.Range("F275").Value = 100112052
$ws.Range("G275").Value = "Albahaca"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 420
$ws.Range("K275").Value = 4500
$ws.Range("L275").Value = 5000
$ws.Range("M275").Value = 4851
$ws.Range("N275").Value = "$/docena de matas"
$ws.Range("O275").Value = "Región Metropolitana"
$ws.Range("P275").Value = 808
$ws.Range("Q275").Value = 6
$ws.Range("R275").Value = "Hortaliza"

# Row 276
$ws.Range("A276").Value = 6
$ws.Range("B276").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C276").Value = "Metropolitana"
$ws.Range("D276").Value = 44179
$ws.Range("E276").Value = 13
$ws.Range("F276").Value = 100112052
$ws.Range("G276").Value = "Albahaca"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Segunda"
$ws.Range("J276").Value = 180
$ws.Range("K276").Value = 4000
$ws.Range("L276").Value = 4000
$ws.Range("M276").Value = 4000
$ws.Range("N276").Value = "$/docena de matas"
$ws.Range("O276").Value = "Región Metropolitana"
$ws.Range("P276").Value = 667
$ws.Range("Q276").Value = 6
$ws.Range("R276").Value = "Hortaliza"
